# Auto-generated edit script: updates numeric cell values per the
# scheduled-runner price/profit refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$cellMap = @{
    "H80" = 496.05
    "I80" = 416.76923
    "J80" = 643.2857
    "K80" = 1250.30769
    "L80" = 1929.8571
    "M80" = -252.3076900000001
    "N80" = -3925.8571
    "H83" = 496.05
    "I83" = 416.76923
    "J83" = 643.2857
    "K83" = 3750.92307
    "L83" = 5789.571300000001
    "M83" = 1241.07693
    "N83" = -15773.5713
    "H112" = 7693185
    "I112" = 674
    "J112" = 8334227.5
    "K112" = 2022
    "L112" = 25002682.5
    "M112" = -914
    "N112" = -25004898.5
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$cellMap = @{
    "H2" = 2723.3333
    "I2" = 3383.1667
    "J2" = 1403.6666
    "K2" = 3383.1667
    "L2" = 1403.6666
    "M2" = -3270.1667
    "N2" = -1629.6666
    "H45" = 50637.332
    "I45" = 75006
    "J45" = 1900
    "K45" = 75006
    "L45" = 1900
    "M45" = -74629
    "N45" = -2654
    "H61" = 2872.889
    "I61" = 2732
    "J61" = 4000
    "K61" = 2732
    "L61" = 4000
    "M61" = -2520
    "N61" = -4424
    "H74" = 1601.6285
    "I74" = 1920.4546
    "J74" = 1455.5
    "K74" = 1920.4546
    "L74" = 1455.5
    "M74" = -1046.4546
    "N74" = -3203.5
    "H77" = 1601.6285
    "I77" = 1920.4546
    "J77" = 1455.5
    "K77" = 9602.273000000001
    "L77" = 7277.5
    "M77" = -5234.273000000001
    "N77" = -16013.5
    "H116" = 2723.3333
    "I116" = 3383.1667
    "J116" = 1403.6666
    "K116" = 3383.1667
    "L116" = 1403.6666
    "M116" = -1089.1667
    "N116" = -5991.6666
    "H136" = 2872.889
    "I136" = 2732
    "J136" = 4000
    "K136" = 8196
    "L136" = 12000
    "M136" = -5646
    "N136" = -17100
    "H139" = 47532.855
    "J139" = 47532.855
    "L139" = 47532.855
    "N139" = -57812.855
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$cellMap = @{
    "H3" = 2723.3333
    "I3" = 3383.1667
    "J3" = 1403.6666
    "K3" = 3383.1667
    "L3" = 1403.6666
    "M3" = -3269.1667
    "N3" = -1631.6666
    "H10" = 28671.818
    "J10" = 29999.047
    "L10" = 29999.047
    "N10" = -30279.047
    "H74" = 52752.855
    "J74" = 52752.855
    "L74" = 52752.855
    "N74" = -54624.855
    "H77" = 52752.855
    "J77" = 52752.855
    "L77" = 158258.565
    "N77" = -167618.565
    "H132" = 34157.5
    "J132" = 34157.5
    "L132" = 34157.5
    "N132" = -44277.5
    "H134" = 3423.8462
    "I134" = 4102
    "K134" = 12306
    "M134" = -9771
    "H138" = 58000.77
    "J138" = 58000.77
    "L138" = 58000.77
    "N138" = -68280.76999999999
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$cellMap = @{
    "H31" = 3127.087
    "I31" = 2034
    "J31" = 3896.2964
    "K31" = 2034
    "L31" = 3896.2964
    "M31" = -1739
    "N31" = -4486.2964
    "H34" = 3127.087
    "I34" = 2034
    "J34" = 3896.2964
    "K34" = 2034
    "L34" = 3896.2964
    "M34" = -1832
    "N34" = -4300.2964
    "H58" = 1259.9
    "I58" = 1122.1111
    "J58" = 2500
    "K58" = 1122.1111
    "L58" = 2500
    "M58" = -919.1111000000001
    "N58" = -2906
    "H132" = 2593.739
    "I132" = 2297.842
    "K132" = 6893.526
    "M132" = -4363.526
    "H134" = 3542.6
    "I134" = 3526.238
    "J134" = 3628.5
    "K134" = 10578.714
    "L134" = 10885.5
    "M134" = -8043.714
    "N134" = -15955.5
    "H136" = 1259.9
    "I136" = 1122.1111
    "J136" = 2500
    "K136" = 3366.3333
    "L136" = 7500
    "M136" = -816.3333000000002
    "N136" = -12600
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$cellMap = @{
    "H68" = 2822.3442
    "I68" = 5160
    "J68" = 1306.027
    "K68" = 15480
    "L68" = 3918.081
    "M68" = -14669
    "N68" = -5540.081
    "H71" = 2822.3442
    "I71" = 5160
    "J71" = 1306.027
    "K71" = 46440
    "L71" = 11754.243
    "M71" = -42384
    "N71" = -19866.243
    "H131" = 25582512
    "J131" = 30304350
    "L131" = 90913050
    "N131" = -90923130
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$cellMap = @{
    "H102" = 2122.6667
    "I102" = 1500
    "J102" = 2200.5
    "K102" = 1500
    "L102" = 2200.5
    "M102" = 122
    "N102" = -5444.5
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$cellMap = @{
    "H16" = 1000
    "I16" = 1000
    "K16" = 1000
    "M16" = -830
    "H125" = 44000
    "J125" = 44000
    "L125" = 44000
    "N125" = -53840
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$cellMap = @{
    "H136" = 1440.7142
    "I136" = 830
    "J136" = 2356.7856
    "K136" = 2490
    "L136" = 7070.3568
    "M136" = 60
    "N136" = -12170.3568
}
foreach ($addr in $cellMap.Keys) {
    $ws.Range($addr).Value = $cellMap[$addr]
}

Write-Output "Updated cells across sheets."